$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 23.82183074951172
$ws.Range("C3").Value = 16.85309410095215
$ws.Range("C4").Value = 16.24298095703125
$ws.Range("C5").Value = 16.1902904510498
$ws.Range("C6").Value = 16.45183563232422
